$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header / metabolite-table preamble (row 6) ---
$ws.Range("A6").Value = "Metabolites reported per g of fresh weight of 6-week-old plant leaf rosettes"
$ws.Range("I6").Value = "C00089"
$ws.Range("J6").Value = "C01793"

# --- Column headers (row 7) ---
$ws.Range("H7").Value = "Starch " + [char]10 + "(mg/g FW)"
$ws.Range("I7").Value = "Sucrose (mg/g FW)"
$ws.Range("J7").Value = "Cholorophyll (mg/g FW)"

# --- SD (short day) measurement rows 8-12 ---
$ws.Range("G8").Value = 0.1206
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 1.2
$ws.Range("J8").Value = 1.8

$ws.Range("G9").Value = 0.1275
$ws.Range("H9").Value = 6.5
$ws.Range("I9").Value = 1.1
$ws.Range("J9").Value = 1.6

$ws.Range("G10").Value = 0.2872
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1.4

$ws.Range("G11").Value = 0.1524
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 0.6
$ws.Range("J11").Value = 2

$ws.Range("G12").Value = 0.2035

# --- LD (long day) measurement rows 13-17: Starch column reverts to "NA" ---
$ws.Range("G13").Value = "NA"
$ws.Range("G14").Value = "NA"
$ws.Range("G15").Value = "NA"
$ws.Range("G16").Value = "NA"
$ws.Range("G17").Value = "NA"

# --- Column widths for the newly populated H:J columns ---
$ws.Columns.Item(8).AutoFit() | Out-Null
$ws.Columns.Item(9).AutoFit() | Out-Null
$ws.Columns.Item(10).AutoFit() | Out-Null

# --- Restore the cursor/selection to K11 ---
$ws.Range("K11").Select()
